# Introduction to unsafe C# - apply "Removed unused elements from slides" commit
#
# 1. Update the datetimeFigureOut placeholder text (6/22/2022 -> 10/8/2022) on
#    the slide master and every slide layout.
# 2. Move/resize the title placeholder on slide 1 (explicit xfrm).
# 3. Split the "Rev. 1 (2022-06-22), " run on slide 1 into "Rev. 2 " and
#    "(2022-10-??), ".
# 4. Remove the unused, empty "Content Placeholder 2" shape from the slides
#    that still had one left over from the layout.

$p = $ppt.ActivePresentation

# --- 1. Date field text on master + all layouts -----------------------------
$newDate = "10/8/2022"

foreach ($sh in $p.SlideMaster.Shapes.Placeholders) {
    if ($sh.PlaceholderFormat.Type -eq 16) {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($i)
    foreach ($sh in $layout.Shapes.Placeholders) {
        if ($sh.PlaceholderFormat.Type -eq 16) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 2. Slide 1 title placeholder gets an explicit position/size ------------
$slide1 = $p.Slides.Item(1)
$title = $slide1.Shapes.Item(1)
# EMU -> points (1 pt = 12700 EMU); nudge by a hair so the float round-trip
# lands on the exact target EMU value instead of one EMU short.
$title.Left = (1154955 / 12700) + 0.00005
$title.Top = (1451020 / 12700) + 0.00005
$title.Width = (8825658 / 12700) + 0.00005
$title.Height = (3329581 / 12700) + 0.00005

# --- 3. Split the revision/date run on slide 1 -------------------------------
$revShape = $slide1.Shapes.Item(4)
$tr = $revShape.TextFrame.TextRange
$revPart = $tr.Characters(1, 7)
$revPart.Text = "Rev. 2 "
$datePart = $tr.Characters(8, 14)
$datePart.Text = "(2022-10-??), "

# --- 4. Remove unused "Content Placeholder 2" shapes -------------------------
$slideNumbers = @(3, 4, 7, 8, 9, 10, 12, 13, 14, 15, 17, 18, 19, 23, 25, 26, 27)
foreach ($n in $slideNumbers) {
    $s = $p.Slides.Item($n)
    $ph = $s.Shapes.Item("Content Placeholder 2")
    $ph.Delete()
}
